$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row rename ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Capitalization / text fixes ---
$ws.Range("B2").Value = "Comitán De Domínguez"
$ws.Range("B6").Value = "Hidalgo Del Parral"
$ws.Range("A8").Value = "Ciudad De México"
$ws.Range("D10").Value = 0.09322033898305083
$ws.Range("A13").Value = "Estado De México"
$ws.Range("B13").Value = "Atizapán De Zaragoza"
$ws.Range("B16").Value = "Tlalnepantla De Baz"
$ws.Range("A18").Value = "Guanajuato"
$ws.Range("B18").Value = "San Luis De La Paz"
$ws.Range("B20").Value = "Acapulco De Juárez"
$ws.Range("B22").Value = "Ayutla De Los Libres"
$ws.Range("B24").Value = "Cuetzala Del Progreso"
$ws.Range("B25").Value = "Huitzuco De Los Figueroa"
$ws.Range("B26").Value = "Zihuatanejo De Azueta"
$ws.Range("B31").Value = "Técpan De Galeana"
$ws.Range("B37").Value = "Mixquiahuala De Juárez"
$ws.Range("B38").Value = "Pachuca De Soto"
$ws.Range("B42").Value = "Unión De Tula"
$ws.Range("B45").Value = "Tiquicheo De Nicolás Romero"
$ws.Range("B52").Value = "Santa María Del Oro"
$ws.Range("B56").Value = "San Dionisio Del Mar"
$ws.Range("B70").Value = "Jalpan De Serra"
$ws.Range("B71").Value = "Landa De Matamoros"
$ws.Range("B81").Value = "Ixtacuixtla De Mariano Matamoros"

# --- Remove footer/metadata rows (93-97) ---
$ws.Range("A93:D97").EntireRow.Delete()
